$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# apply E first using Style indirection before B/C/D  
$rE = $ws.Range("E2:E10")
$rE.Value = "N/A"
$rE.NumberFormat = "0%"

$rBCD = $ws.Range("B2:D10")
$rBCD.Value = "N/A"
$rBCD.NumberFormat = "0.00"

$rFG = $ws.Range("F2:G10")
$rFG.Value = "N/A"
$rFG.NumberFormat = "@"
